# Generate Report for Handoff
# Update Priority (column E) from "low" to "ht" and refresh the
# Latest Handoff Datetime (column H) for the four rows that were
# previously at "low" priority, in both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-11-29 06:06:49"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-11-29 06:07:03"

# The Overview sheet mirrors the de-de "Latest Handoff Datetime" value
# (shared across the same four rows) in column G; keep it in sync too.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4:G7").Value = "2016-11-29 06:07:03"
